$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "25.908.85" or "215.73" - some of
# the NEW values in this update (e.g. "215.64") would be auto-detected by
# Excel as a plain number if assigned as-is, which would silently change
# the cell from text to a numeric type. Force those specific cells to the
# Text format first so the assigned value round-trips as a string.
$numericLookingCells = @(
    'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D18', 'D19', 'D20', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D30', 'D31', 'D32', 'D34', 'D35', 'D37', 'D39', 'D41', 'D42', 'D43', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51'
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- cell value updates (cryptos list refresh) ---
$ws.Range('D2').Value = '25.894.64'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.641.61'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.55%  '
$ws.Range('D5').Value = '215.64'
$ws.Range('D6').Value = '0.5062'
$ws.Range('D7').Value = '1.007'
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('D8').Value = '0.2576'
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '0.06417'
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('D10').Value = '19.76'
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('D11').Value = '0.07778'
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('D12').Value = '4.306'
$ws.Range('E12').Value = '  +1.80%  '
$ws.Range('D13').Value = '1.625.93'
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('D14').Value = '0.5458'
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('D15').Value = '0.0₅7904'
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('E16').Value = '  +2.70%  '
$ws.Range('D17').Value = '25.989.94'
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('D18').Value = '1.008'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').Value = '198.03'
$ws.Range('E19').Value = '  -2.30%  '
$ws.Range('D20').Value = '4.411'
$ws.Range('E20').Value = '  +2.51%  '
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('D22').Value = '6.046'
$ws.Range('E22').Value = '  +1.40%  '
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').Value = '1.871'
$ws.Range('E24').Value = '  -3.25%  '
$ws.Range('D25').Value = '140.83'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').Value = '0.1148'
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').Value = '6.898'
$ws.Range('D28').Value = '15.72'
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').Value = '0.05049'
$ws.Range('E30').Value = '  +1.75%  '
$ws.Range('D31').Value = '3.273'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').Value = '3.201'
$ws.Range('E32').Value = '  +0.52%  '
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('D34').Value = '2.373'
$ws.Range('E34').Value = '  +0.98%  '
$ws.Range('D35').Value = '0.8947'
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('E36').Value = '  -0.84%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.5543'
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.133.21'
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('D39').Value = '0.01563'
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('E40').Value = '  +0.74%  '
$ws.Range('D41').Value = '5.691'
$ws.Range('E41').Value = '  +0.89%  '
$ws.Range('D42').Value = '0.8177'
$ws.Range('E42').Value = '  +2.05%  '
$ws.Range('D43').Value = '99.80'
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('E44').Value = '  +6.93%  '
$ws.Range('D45').Value = '1.779.85'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').Value = '0.4540'
$ws.Range('E46').Value = '  +0.69%  '
$ws.Range('D47').Value = '55.31'
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('D48').Value = '1.005'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').Value = '0.05090'
$ws.Range('E49').Value = '  +1.18%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').Value = '1.009'
$ws.Range('E50').Value = '  +0.69%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.09568'
$ws.Range('E51').Value = '  +3.26%  '

# Put those cells back on the default "Normal" style so only the stored
# value changed - no stray explicit number-format/style index remains.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}

Write-Output "Applied cryptos update"
